# Prep for LS1 submission page using broker1
#
# - Rename "aclistingdata1" -> "aclistingdata"
# - Add a new, empty worksheet "brokerls1data" right after it
# - Make the new sheet the active tab, with its own selection
# - Leave the old active sheet's selection parked at E34 (no longer the
#   selected tab)

$wb = $excel.ActiveWorkbook

# The workbook has two sheets: "credentials" (1) and "aclistingdata1" (2).
$acSheet = $wb.Worksheets.Item(2)
$acSheet.Name = "aclistingdata"

# Insert the new sheet immediately after "aclistingdata".
$brokerSheet = $wb.Worksheets.Add($null, $acSheet)
$brokerSheet.Name = "brokerls1data"

# Park the selection on the (now inactive) aclistingdata sheet.
$acSheet.Activate()
$null = $acSheet.Range("E34").Select()

# Make the new broker sheet the active tab with its own selection.
$brokerSheet.Activate()
$null = $brokerSheet.Range("F27").Select()
